$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.9
$ws.Range("P3").Value = 1.44
$ws.Range("Q3").Value = 2.63
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.83
$ws.Range("T3").Value = 7
$ws.Range("V3").Value = 9.5
$ws.Range("W3").Value = 21
$ws.Range("Z3").Value = 9
$ws.Range("AG3").Value = 34
$ws.Range("AH3").Value = 26

# Row 4
$ws.Range("G4").Value = 2.45
$ws.Range("I4").Value = 2.8
$ws.Range("T4").Value = 7.5
$ws.Range("W4").Value = 23
$ws.Range("AF4").Value = 11

# Row 5
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 1.9
$ws.Range("J5").Value = 1.14
$ws.Range("K5").Value = 5.5
$ws.Range("L5").Value = 1.67
$ws.Range("M5").Value = 2.2
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 1.67
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 2.63
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 8.5
$ws.Range("V5").Value = 19
$ws.Range("Y5").Value = 67
$ws.Range("Z5").Value = 5.5
$ws.Range("AD5").Value = 4.75
$ws.Range("AF5").Value = 10

# Row 6
$ws.Range("G6").Value = 2.15
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 1.06
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 3.4
$ws.Range("N6").Value = 2.03
$ws.Range("O6").Value = 1.83
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.91
$ws.Range("U6").Value = 10
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 10
$ws.Range("AC6").Value = 51
$ws.Range("AE6").Value = 17
$ws.Range("AH6").Value = 26
$ws.Range("AI6").Value = 34
$ws.Range("AJ6").Value = 251

# Row 7
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 2.9
$ws.Range("J7").Value = 1.1
$ws.Range("K7").Value = 7
$ws.Range("L7").Value = 1.5
$ws.Range("M7").Value = 2.63
$ws.Range("N7").Value = 2.5
$ws.Range("O7").Value = 1.53
$ws.Range("P7").Value = 1.57
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 2.1
$ws.Range("S7").Value = 1.67
$ws.Range("T7").Value = 6.5
$ws.Range("U7").Value = 11
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 41
$ws.Range("Z7").Value = 7
$ws.Range("AA7").Value = 6
$ws.Range("AB7").Value = 19
$ws.Range("AC7").Value = 67
$ws.Range("AD7").Value = 7
$ws.Range("AF7").Value = 12
$ws.Range("AG7").Value = 34
$ws.Range("AH7").Value = 29
$ws.Range("AI7").Value = 41
$ws.Range("AJ7").Value = 501

# Row 8
$ws.Range("H8").Value = 3.9
$ws.Range("J8").Value = 1.07
$ws.Range("K8").Value = 9
$ws.Range("L8").Value = 1.4
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 2.2
$ws.Range("O8").Value = 1.67
$ws.Range("R8").Value = 2.2
$ws.Range("S8").Value = 1.62
$ws.Range("T8").Value = 5.5
$ws.Range("Z8").Value = 9
$ws.Range("AC8").Value = 81

# Row 9
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 3.2
$ws.Range("J9").Value = 1.08
$ws.Range("K9").Value = 8
$ws.Range("X9").Value = 23
$ws.Range("Y9").Value = 41
$ws.Range("AA9").Value = 6
$ws.Range("AG9").Value = 29
$ws.Range("AH9").Value = 29

# Row 10
$ws.Range("G10").Value = 2.63
$ws.Range("H10").Value = 3.25
$ws.Range("I10").Value = 2.63
$ws.Range("U10").Value = 12
$ws.Range("W10").Value = 26
$ws.Range("X10").Value = 23
$ws.Range("AB10").Value = 15
$ws.Range("AD10").Value = 7.5
$ws.Range("AE10").Value = 12
$ws.Range("AF10").Value = 10
$ws.Range("AG10").Value = 26

# Row 11
$ws.Range("G11").Value = 4.2
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 1.75
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 13
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4.33
$ws.Range("N11").Value = 1.67
$ws.Range("O11").Value = 2.2
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 2.2
$ws.Range("U11").Value = 23
$ws.Range("V11").Value = 13
$ws.Range("W11").Value = 41
$ws.Range("X11").Value = 29
$ws.Range("AA11").Value = 7.5
$ws.Range("AE11").Value = 9.5
$ws.Range("AG11").Value = 15
$ws.Range("AH11").Value = 13

# Row 12
$ws.Range("G12").Value = 1.65
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 5
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 1.75
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 7
$ws.Range("W12").Value = 12
$ws.Range("AA12").Value = 7.5
$ws.Range("AB12").Value = 19
$ws.Range("AC12").Value = 67
$ws.Range("AD12").Value = 12
$ws.Range("AE12").Value = 26
$ws.Range("AF12").Value = 17
$ws.Range("AJ12").Value = 401

# Row 13
$ws.Range("L13").Value = 1.57
$ws.Range("M13").Value = 2.25

# Row 14
$ws.Range("G14").Value = 3.1
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.3
$ws.Range("J14").Value = 1.08
$ws.Range("K14").Value = 7.5
$ws.Range("L14").Value = 1.4
$ws.Range("M14").Value = 2.75
$ws.Range("N14").Value = 2.35
$ws.Range("O14").Value = 1.57
$ws.Range("P14").Value = 1.53
$ws.Range("Q14").Value = 2.38
$ws.Range("R14").Value = 2
$ws.Range("S14").Value = 1.75
$ws.Range("T14").Value = 8
$ws.Range("Z14").Value = 7.5
$ws.Range("AA14").Value = 6
$ws.Range("AC14").Value = 67
$ws.Range("AD14").Value = 6.5
$ws.Range("AH14").Value = 21
$ws.Range("AJ14").Value = 451

# Row 15
$ws.Range("N15").Value = 1.98
$ws.Range("O15").Value = 1.88

# Row 16
$ws.Range("H16").Value = 4.2
$ws.Range("K16").Value = 15
$ws.Range("L16").Value = 1.2
$ws.Range("M16").Value = 4.33
$ws.Range("N16").Value = 1.65
$ws.Range("O16").Value = 2.2
$ws.Range("P16").Value = 1.3
$ws.Range("Q16").Value = 3.4
$ws.Range("R16").Value = 1.7
$ws.Range("S16").Value = 2.05
$ws.Range("T16").Value = 15
$ws.Range("Z16").Value = 15
$ws.Range("AD16").Value = 8.5
$ws.Range("AE16").Value = 8.5
$ws.Range("AJ16").Value = 201

# Row 17
$ws.Range("K17").Value = 9

# Row 18
$ws.Range("G18").Value = 1.91
$ws.Range("I18").Value = 3.7
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("N18").Value = 2.08
$ws.Range("O18").Value = 1.73
$ws.Range("R18").Value = 1.91
$ws.Range("S18").Value = 1.8
$ws.Range("U18").Value = 8.5
$ws.Range("W18").Value = 17
$ws.Range("Z18").Value = 9
$ws.Range("AC18").Value = 51
$ws.Range("AD18").Value = 9.5
$ws.Range("AE18").Value = 19
$ws.Range("AF18").Value = 13
$ws.Range("AH18").Value = 34
$ws.Range("AJ18").Value = 351

# Row 32
$ws.Range("G32").Value = 1.42
$ws.Range("I32").Value = 8
$ws.Range("R32").Value = 2.25
$ws.Range("S32").Value = 1.57
$ws.Range("T32").Value = 5.5
$ws.Range("AB32").Value = 23
$ws.Range("AG32").Value = 101
$ws.Range("AH32").Value = 67
